# ADUserdata.xlsx — "Add files via upload" edit
#
# Renames the "投放开始日期/投放结束日期" (start/end date) headers to the
# shorter "投放日期/结束日期", re-types a couple of the date values in the
# sample data rows (row 2 and row 38), and lets the column widths re-flow
# to match the new (narrower) header/content widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (both the left A:I block and the mirrored J:Q block) ---
$ws.Range("G1").Value = "投放日期"
$ws.Range("H1").Value = "结束日期"
$ws.Range("P1").Value = "投放日期"
$ws.Range("Q1").Value = "结束日期"

# --- Row 2 sample data: P2 becomes a literal (unparsable) text date, ---
# --- Q2 moves two months later                                      ---
$ws.Range("P2").Value = "2026/13/26"
$ws.Range("Q2").Value = "12/09/2099"

# --- Row 38 sample data: both dates shift forward ---
$ws.Range("G38").Value = "10/26/2026"
$ws.Range("H38").Value = "12/09/2020"

# --- Column widths re-flow for the shortened headers / retyped dates ---
$ws.Range("E1").ColumnWidth = 12.4
$ws.Range("G1").ColumnWidth = 10.86
$ws.Range("H1").ColumnWidth = 10.86
$ws.Range("N1").ColumnWidth = 12.4
$ws.Range("P1").ColumnWidth = 10.86
$ws.Range("Q1").ColumnWidth = 10.86

# --- Clear the stale selection left over from editing, back to A1 ---
[void]$ws.Range("A1").Select()
